$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...by contacting me at any time before the 20th of July 2020." ->
#    "...by contacting me at any time before up to 1 week after the
#     interview has taken place." (in the main "What are your rights" body
#    paragraph, not the bullet list further down which already reads this
#    way).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " by contacting me at any time before the 20th of July 2020",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " up to 1 week after the interview has taken place", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2,3,4) "...receive a copy of the interview transcript which will be
#    emailed to you, and you will have the opportunity to add comments to
#    the interview up to 2 weeks after receiving the transcript." ->
#    "...receive a copy of the interview recording which will be emailed to
#    you, and you will have the opportunity to edit, appended and remove
#    details up to 1 week after receiving the recording."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "choose to receive a copy of the interview transcript which will be emailed to you, and you will have the opportunity to add comments to the interview up to 2 weeks after receiving the transcript",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "choose to receive a copy of the interview recording which will be emailed to you, and you will have the opportunity to edit, appended and remove details up to 1 week after receiving the recording",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5) "...when a full transcript of your interview is sent..." ->
#    "...when a recording of your interview is sent..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "full transcript of",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "recording of", 2) | Out-Null

# ---------------------------------------------------------------------------
# 6) "I will read the notes or transcript of the interview. The interview
#    transcripts, summaries and any recordings will be kept securely and
#    destroyed" ->
#    "I will read the notes of the interview. The interview summaries and
#    any recordings will be kept securely and destroyed"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "I will read the notes or transcript of the interview. The interview transcripts, summaries and any recordings will be kept securely and destroyed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "I will read the notes of the interview. The interview summaries and any recordings will be kept securely and destroyed",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 7) Bullet list: drop the "receive a copy of your interview transcript;"
#    bullet entirely (it is now redundant with the "...interview recording;"
#    bullet above it) and merge its paragraph mark into the following
#    "edit/append/remove..." bullet, which also gains "up to 1 week " before
#    "after the interview."
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "receive a copy of your interview transcript;") {
        $r = $d.Range($para.Range.Start, $para.Range.End)
        $r.Delete()
        $found = $true
        break
    }
}
Write-Output "Removed transcript bullet: $found"

$d.Content.Find.Execute(
    "edit/append/remove any details after the interview.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "edit/append/remove any details up to 1 week after the interview.",
    2) | Out-Null

Write-Output "done"
